# Power_BusInfo.xlsx edit
# Commit: "Fix pyomo vLineP, vLineQ, vGenQ handeling. Add bounds for quadratic variables"
#
# The susceptance (Bs, column H) and conductance (Gs, column I) bound values
# for every bus row (rows 8-16 on the "ScenarioA" sheet) were changed from 0
# to 1, giving the quadratic variables an actual (non-zero) bound instead of
# the previous placeholder of 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioA")

# Add bounds for the quadratic variables: Bs (H) and Gs (I) go from 0 -> 1
# for all bus data rows (8 through 16).
$ws.Range("H8:I16").Value = 1

# Leave the cursor where the author left it before saving.
$ws.Range("H22").Select() | Out-Null
